{"js": "const pairs = [\n  [\"85-7=\", \"59-21=\"],\n  [\"62+4=\", \"58-4=\"],\n  [\"82-67=\", \"69-22=\"],\n  [\"88+1=\", \"97-90=\"],\n  [\"82-68=\", \"11+14=\"],\n  [\"76-22=\", \"16-15=\"],\n  [\"67+16=\", \"11+52=\"],\n  [\"2+75=\", \"49+39=\"],\n  [\"78-34=\", \"28+58=\"],\n  [\"35+0=\", \"47+8=\"],\n  [\"46+6=\", \"56-47=\"],\n  [\"31+49=\", \"52+45=\"],\n  [\"40-18=\", \"46-37=\"],\n  [\"70-39=\", \"65+2=\"],\n  [\"61-58=\", \"80-40=\"],\n  [\"8+51=\", \"33+53=\"],\n  [\"87-1=\", \"79-72=\"],\n  [\"21+61=\", \"72-50=\"],\n  [\"70-4=\", \"21-5=\"],\n  [\"3+18=\", \"19+80=\"],\n  [\"39+58=\", \"5+32=\"],\n  [\"51-20=\", \"76-8=\"],\n  [\"71+16=\", \"79-75=\"],\n  [\"45+29=\", \"99-66=\"],\n  [\"69-52=\", \"20+77=\"],\n  [\"41+37=\", \"17+7=\"],\n  [\"16+68=\", \"23+59=\"],\n  [\"22+7=\", \"78-16=\"],\n  [\"27-23=\", \"39-17=\"],\n  [\"31+51=\", \"85-35=\"],\n  [\"75-72=\", \"6+61=\"],\n  [\"80+15=\", \"48-11=\"],\n  [\"98-48=\", \"25+26=\"],\n  [\"43+34=\", \"4+39=\"],\n  [\"31+26=\", \"36-35=\"],\n  [\"98-56=\", \"75-12=\"],\n  [\"90-28=\", \"66+31=\"],\n  [\"75-47=\", \"70-33=\"],\n  [\"42-10=\", \"74-67=\"],\n  [\"50+2=\", \"59-25=\"],\n  [\"59-18=\", \"62-8=\"],\n  [\"22+10=\", \"26+38=\"],\n  [\"95-14=\", \"59-7=\"],\n  [\"47+45=\", \"76-36=\"],\n  [\"4+85=\", \"19+48=\"],\n  [\"65-59=\", \"0+29=\"],\n  [\"73-7=\", \"89+3=\"],\n  [\"24-22=\", \"55-47=\"],\n  [\"32+63=\", \"83-9=\"],\n  [\"79+9=\", \"1+27=\"],\n  [\"4+77=\", \"35+2=\"],\n  [\"18-18=\", \"75+22=\"],\n  [\"32+58=\", \"32+66=\"],\n  [\"23+7=\", \"88-42=\"],\n  [\"42+42=\", \"85-22=\"],\n  [\"51+2=\", \"81-70=\"],\n  [\"93-35=\", \"85+12=\"],\n  [\"86-19=\", \"53-42=\"],\n  [\"22+59=\", \"56+27=\"],\n  [\"77+7=\", \"18-13=\"],\n  [\"28+15=\", \"99-96=\"],\n  [\"24-13=\", \"87-64=\"],\n  [\"19+58=\", \"84+0=\"],\n  [\"26+3=\", \"56-19=\"],\n  [\"94-54=\", \"43-2=\"],\n  [\"81+13=\", \"69+25=\"],\n  [\"36+38=\", \"96-69=\"],\n  [\"44-28=\", \"14-5=\"],\n  [\"34+58=\", \"40+33=\"],\n  [\"94-15=\", \"40+56=\"],\n  [\"54+41=\", \"18+65=\"],\n  [\"17+4=\", \"80-14=\"],\n  [\"79+17=\", \"49-24=\"],\n  [\"13+59=\", \"76-60=\"],\n  [\"79-53=\", \"74-22=\"],\n  [\"92-71=\", \"96-53=\"],\n  [\"91-2=\", \"60-37=\"],\n  [\"64-1=\", \"67-63=\"],\n  [\"58-46=\", \"2+94=\"],\n  [\"25+41=\", \"40-5=\"],\n  [\"69-8=\", \"23+23=\"],\n  [\"94+5=\", \"82-19=\"],\n  [\"45-44=\", \"69-56=\"],\n  [\"21-4=\", \"11+50=\"],\n  [\"6+39=\", \"83+2=\"],\n  [\"34+46=\", \"56+8=\"],\n  [\"22+11=\", \"42-0=\"],\n  [\"38+46=\", \"45-24=\"],\n  [\"74-18=\", \"77-50=\"],\n  [\"14+36=\", \"21-13=\"],\n  [\"42+55=\", \"66-62=\"],\n  [\"66-29=\", \"49+22=\"],\n  [\"95-80=\", \"21+55=\"],\n  [\"40+7=\", \"78+1=\"],\n  [\"48-23=\", \"70-58=\"],\n  [\"51-14=\", \"31+33=\"],\n  [\"43+3=\", \"35+10=\"],\n  [\"95-50=\", \"57+12=\"],\n  [\"87-5=\", \"42-38=\"],\n  [\"58-9=\", \"26+31=\"]\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$d.Content.Find.Execute(\"85-7=\", $false, $false, $false, $false, $false, $true, 1, $false, \"59-21=\", 2) | Out-Null\n$d.Content.Find.Execute(\"62+4=\", $false, $false, $false, $false, $false, $true, 1, $false, \"58-4=\", 2) | Out-Null\n$d.Content.Find.Execute(\"82-67=\", $false, $false, $false, $false, $false, $true, 1, $false, \"69-22=\", 2) | Out-Null\n$d.Content.Find.Execute(\"88+1=\", $false, $false, $false, $false, $false, $true, 1, $false, \"97-90=\", 2) | Out-Null\n$d.Content.Find.Execute(\"82-68=\", $false, $false, $false, $false, $false, $true, 1, $false, \"11+14=\", 2) | Out-Null\n$d.Content.Find.Execute(\"76-22=\", $false, $false, $false, $false, $false, $true, 1, $false, \"16-15=\", 2) | Out-Null\n$d.Content.Find.Execute(\"67+16=\", $false, $false, $false, $false, $false, $true, 1, $false, \"11+52=\", 2) | Out-Null\n$d.Content.Find.Execute(\"2+75=\", $false, $false, $false, $false, $false, $true, 1, $false, \"49+39=\", 2) | Out-Null\n$d.Content.Find.Execute(\"78-34=\", $false, $false, $false, $false, $false, $true, 1, $false, \"28+58=\", 2) | Out-Null\n$d.Content.Find.Execute(\"35+0=\", $false, $false, $false, $false, $false, $true, 1, $false, \"47+8=\", 2) | Out-Null\n$d.Content.Find.Execute(\"46+6=\", $false, $false, $false, $false, $false, $true, 1, $false, \"56-47=\", 2) | Out-Null\n$d.Content.Find.Execute(\"31+49=\", $false, $false, $false, $false, $false, $true, 1, $false, \"52+45=\", 2) | Out-Null\n$d.Content.Find.Execute(\"40-18=\", $false, $false, $false, $false, $false, $true, 1, $false, \"46-37=\", 2) | Out-Null\n$d.Content.Find.Execute(\"70-39=\", $false, $false, $false, $false, $false, $true, 1, $false, \"65+2=\", 2) | Out-Null\n$d.Content.Find.Execute(\"61-58=\", $false, $false, $false, $false, $false, $true, 1, $false, \"80-40=\", 2) | Out-Null\n$d.Content.Find.Execute(\"8+51=\", $false, $false, $false, $false, $false, $true, 1, $false, \"33+53=\", 2) | Out-Null\n$d.Content.Find.Execute(\"87-1=\", $false, $false, $false, $false, $false, $true, 1, $false, \"79-72=\", 2) | Out-Null\n$d.Content.Find.Execute(\"21+61=\", $false, $false, $false, $false, $false, $true, 1, $false, \"72-50=\", 2) | Out-Null\n$d.Content.Find.Execute(\"70-4=\", $false, $false, $false, $false, $false, $true, 1, $false, \"21-5=\", 2) | Out-Null\n$d.Content.Find.Execute(\"3+18=\", $false, $false, $false, $false, $false, $true, 1, $false, \"19+80=\", 2) | Out-Null\n$d.Content.Find.Execute(\"39+58=\", $false, $false, $false, $false, $false, $true, 1, $false, \"5+32=\", 2) | Out-Null\n$d.Content.Find.Execute(\"51-20=\", $false, $false, $false, $false, $false, $true, 1, $false, \"76-8=\", 2) | Out-Null\n$d.Content.Find.Execute(\"71+16=\", $false, $false, $false, $false, $false, $true, 1, $false, \"79-75=\", 2) | Out-Null\n$d.Content.Find.Execute(\"45+29=\", $false, $false, $false, $false, $false, $true, 1, $false, \"99-66=\", 2) | Out-Null\n$d.Content.Find.Execute(\"69-52=\", $false, $false, $false, $false, $false, $true, 1, $false, \"20+77=\", 2) | Out-Null\n$d.Content.Find.Execute(\"41+37=\", $false, $false, $false, $false, $false, $true, 1, $false, \"17+7=\", 2) | Out-Null\n$d.Content.Find.Execute(\"16+68=\", $false, $false, $false, $false, $false, $true, 1, $false, \"23+59=\", 2) | Out-Null\n$d.Content.Find.Execute(\"22+7=\", $false, $false, $false, $false, $false, $true, 1, $false, \"78-16=\", 2) | Out-Null\n$d.Content.Find.Execute(\"27-23=\", $false, $false, $false, $false, $false, $true, 1, $false, \"39-17=\", 2) | Out-Null\n$d.Content.Find.Execute(\"31+51=\", $false, $false, $false, $false, $false, $true, 1, $false, \"85-35=\", 2) | Out-Null\n$d.Content.Find.Execute(\"75-72=\", $false, $false, $false, $false, $false, $true, 1, $false, \"6+61=\", 2) | Out-Null\n$d.Content.Find.Execute(\"80+15=\", $false, $false, $false, $false, $false, $true, 1, $false, \"48-11=\", 2) | Out-Null\n$d.Content.Find.Execute(\"98-48=\", $false, $false, $false, $false, $false, $true, 1, $false, \"25+26=\", 2) | Out-Null\n$d.Content.Find.Execute(\"43+34=\", $false, $false, $false, $false, $false, $true, 1, $false, \"4+39=\", 2) | Out-Null\n$d.Content.Find.Execute(\"31+26=\", $false, $false, $false, $false, $false, $true, 1, $false, \"36-35=\", 2) | Out-Null\n$d.Content.Find.Execute(\"98-56=\", $false, $false, $false, $false, $false, $true, 1, $false, \"75-12=\", 2) | Out-Null\n$d.Content.Find.Execute(\"90-28=\", $false, $false, $false, $false, $false, $true, 1, $false, \"66+31=\", 2) | Out-Null\n$d.Content.Find.Execute(\"75-47=\", $false, $false, $false, $false, $false, $true, 1, $false, \"70-33=\", 2) | Out-Null\n$d.Content.Find.Execute(\"42-10=\", $false, $false, $false, $false, $false, $true, 1, $false, \"74-67=\", 2) | Out-Null\n$d.Content.Find.Execute(\"50+2=\", $false, $false, $false, $false, $false, $true, 1, $false, \"59-25=\", 2) | Out-Null\n$d.Content.Find.Execute(\"59-18=\", $false, $false, $false, $false, $false, $true, 1, $false, \"62-8=\", 2) | Out-Null\n$d.Content.Find.Execute(\"22+10=\", $false, $false, $false, $false, $false, $true, 1, $false, \"26+38=\", 2) | Out-Null\n$d.Content.Find.Execute(\"95-14=\", $false, $false, $false, $false, $false, $true, 1, $false, \"59-7=\", 2) | Out-Null\n$d.Content.Find.Execute(\"47+45=\", $false, $false, $false, $false, $false, $true, 1, $false, \"76-36=\", 2) | Out-Null\n$d.Content.Find.Execute(\"4+85=\", $false, $false, $false, $false, $false, $true, 1, $false, \"19+48=\", 2) | Out-Null\n$d.Content.Find.Execute(\"65-59=\", $false, $false, $false, $false, $false, $true, 1, $false, \"0+29=\", 2) | Out-Null\n$d.Content.Find.Execute(\"73-7=\", $false, $false, $false, $false, $false, $true, 1, $false, \"89+3=\", 2) | Out-Null\n$d.Content.Find.Execute(\"24-22=\", $false, $false, $false, $false, $false, $true, 1, $false, \"55-47=\", 2) | Out-Null\n$d.Content.Find.Execute(\"32+63=\", $false, $false, $false, $false, $false, $true, 1, $false, \"83-9=\", 2) | Out-Null\n$d.Content.Find.Execute(\"79+9=\", $false, $false, $false, $false, $false, $true, 1, $false, \"1+27=\", 2) | Out-Null\n$d.Content.Find.Execute(\"4+77=\", $false, $false, $false, $false, $false, $true, 1, $false, \"35+2=\", 2) | Out-Null\n$d.Content.Find.Execute(\"18-18=\", $false, $false, $false, $false, $false, $true, 1, $false, \"75+22=\", 2) | Out-Null\n$d.Content.Find.Execute(\"32+58=\", $false, $false, $false, $false, $false, $true, 1, $false, \"32+66=\", 2) | Out-Null\n$d.Content.Find.Execute(\"23+7=\", $false, $false, $false, $false, $false, $true, 1, $false, \"88-42=\", 2) | Out-Null\n$d.Content.Find.Execute(\"42+42=\", $false, $false, $false, $false, $false, $true, 1, $false, \"85-22=\", 2) | Out-Null\n$d.Content.Find.Execute(\"51+2=\", $false, $false, $false, $false, $false, $true, 1, $false, \"81-70=\", 2) | Out-Null\n$d.Content.Find.Execute(\"93-35=\", $false, $false, $false, $false, $false, $true, 1, $false, \"85+12=\", 2) | Out-Null\n$d.Content.Find.Execute(\"86-19=\", $false, $false, $false, $false, $false, $true, 1, $false, \"53-42=\", 2) | Out-Null\n$d.Content.Find.Execute(\"22+59=\", $false, $false, $false, $false, $false, $true, 1, $false, \"56+27=\", 2) | Out-Null\n$d.Content.Find.Execute(\"77+7=\", $false, $false, $false, $false, $false, $true, 1, $false, \"18-13=\", 2) | Out-Null\n$d.Content.Find.Execute(\"28+15=\", $false, $false, $false, $false, $false, $true, 1, $false, \"99-96=\", 2) | Out-Null\n$d.Content.Find.Execute(\"24-13=\", $false, $false, $false, $false, $false, $true, 1, $false, \"87-64=\", 2) | Out-Null\n$d.Content.Find.Execute(\"19+58=\", $false, $false, $false, $false, $false, $true, 1, $false, \"84+0=\", 2) | Out-Null\n$d.Content.Find.Execute(\"26+3=\", $false, $false, $false, $false, $false, $true, 1, $false, \"56-19=\", 2) | Out-Null\n$d.Content.Find.Execute(\"94-54=\", $false, $false, $false, $false, $false, $true, 1, $false, \"43-2=\", 2) | Out-Null\n$d.Content.Find.Execute(\"81+13=\", $false, $false, $false, $false, $false, $true, 1, $false, \"69+25=\", 2) | Out-Null\n$d.Content.Find.Execute(\"36+38=\", $false, $false, $false, $false, $false, $true, 1, $false, \"96-69=\", 2) | Out-Null\n$d.Content.Find.Execute(\"44-28=\", $false, $false, $false, $false, $false, $true, 1, $false, \"14-5=\", 2) | Out-Null\n$d.Content.Find.Execute(\"34+58=\", $false, $false, $false, $false, $false, $true, 1, $false, \"40+33=\", 2) | Out-Null\n$d.Content.Find.Execute(\"94-15=\", $false, $false, $false, $false, $false, $true, 1, $false, \"40+56=\", 2) | Out-Null\n$d.Content.Find.Execute(\"54+41=\", $false, $false, $false, $false, $false, $true, 1, $false, \"18+65=\", 2) | Out-Null\n$d.Content.Find.Execute(\"17+4=\", $false, $false, $false, $false, $false, $true, 1, $false, \"80-14=\", 2) | Out-Null\n$d.Content.Find.Execute(\"79+17=\", $false, $false, $false, $false, $false, $true, 1, $false, \"49-24=\", 2) | Out-Null\n$d.Content.Find.Execute(\"13+59=\", $false, $false, $false, $false, $false, $true, 1, $false, \"76-60=\", 2) | Out-Null\n$d.Content.Find.Execute(\"79-53=\", $false, $false, $false, $false, $false, $true, 1, $false, \"74-22=\", 2) | Out-Null\n$d.Content.Find.Execute(\"92-71=\", $false, $false, $false, $false, $false, $true, 1, $false, \"96-53=\", 2) | Out-Null\n$d.Content.Find.Execute(\"91-2=\", $false, $false, $false, $false, $false, $true, 1, $false, \"60-37=\", 2) | Out-Null\n$d.Content.Find.Execute(\"64-1=\", $false, $false, $false, $false, $false, $true, 1, $false, \"67-63=\", 2) | Out-Null\n$d.Content.Find.Execute(\"58-46=\", $false, $false, $false, $false, $false, $true, 1, $false, \"2+94=\", 2) | Out-Null\n$d.Content.Find.Execute(\"25+41=\", $false, $false, $false, $false, $false, $true, 1, $false, \"40-5=\", 2) | Out-Null\n$d.Content.Find.Execute(\"69-8=\", $false, $false, $false, $false, $false, $true, 1, $false, \"23+23=\", 2) | Out-Null\n$d.Content.Find.Execute(\"94+5=\", $false, $false, $false, $false, $false, $true, 1, $false, \"82-19=\", 2) | Out-Null\n$d.Content.Find.Execute(\"45-44=\", $false, $false, $false, $false, $false, $true, 1, $false, \"69-56=\", 2) | Out-Null\n$d.Content.Find.Execute(\"21-4=\", $false, $false, $false, $false, $false, $true, 1, $false, \"11+50=\", 2) | Out-Null\n$d.Content.Find.Execute(\"6+39=\", $false, $false, $false, $false, $false, $true, 1, $false, \"83+2=\", 2) | Out-Null\n$d.Content.Find.Execute(\"34+46=\", $false, $false, $false, $false, $false, $true, 1, $false, \"56+8=\", 2) | Out-Null\n$d.Content.Find.Execute(\"22+11=\", $false, $false, $false, $false, $false, $true, 1, $false, \"42-0=\", 2) | Out-Null\n$d.Content.Find.Execute(\"38+46=\", $false, $false, $false, $false, $false, $true, 1, $false, \"45-24=\", 2) | Out-Null\n$d.Content.Find.Execute(\"74-18=\", $false, $false, $false, $false, $false, $true, 1, $false, \"77-50=\", 2) | Out-Null\n$d.Content.Find.Execute(\"14+36=\", $false, $false, $false, $false, $false, $true, 1, $false, \"21-13=\", 2) | Out-Null\n$d.Content.Find.Execute(\"42+55=\", $false, $false, $false, $false, $false, $true, 1, $false, \"66-62=\", 2) | Out-Null\n$d.Content.Find.Execute(\"66-29=\", $false, $false, $false, $false, $false, $true, 1, $false, \"49+22=\", 2) | Out-Null\n$d.Content.Find.Execute(\"95-80=\", $false, $false, $false, $false, $false, $true, 1, $false, \"21+55=\", 2) | Out-Null\n$d.Content.Find.Execute(\"40+7=\", $false, $false, $false, $false, $false, $true, 1, $false, \"78+1=\", 2) | Out-Null\n$d.Content.Find.Execute(\"48-23=\", $false, $false, $false, $false, $false, $true, 1, $false, \"70-58=\", 2) | Out-Null\n$d.Content.Find.Execute(\"51-14=\", $false, $false, $false, $false, $false, $true, 1, $false, \"31+33=\", 2) | Out-Null\n$d.Content.Find.Execute(\"43+3=\", $false, $false, $false, $false, $false, $true, 1, $false, \"35+10=\", 2) | Out-Null\n$d.Content.Find.Execute(\"95-50=\", $false, $false, $false, $false, $false, $true, 1, $false, \"57+12=\", 2) | Out-Null\n$d.Content.Find.Execute(\"87-5=\", $false, $false, $false, $false, $false, $true, 1, $false, \"42-38=\", 2) | Out-Null\n$d.Content.Find.Execute(\"58-9=\", $false, $false, $false, $false, $false, $true, 1, $false, \"26+31=\", 2) | Out-Null\n"}
